$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G5"  = 2.3
    "H5"  = 2.9
    "I5"  = 3.5
    "J5"  = 3.2
    "K5"  = 1.83
    "L5"  = 4.5
    "M5"  = 1.13
    "N5"  = 6
    "O5"  = 1.57
    "P5"  = 2.25
    "Q5"  = 2.88
    "R5"  = 1.4
    "S5"  = 6
    "T5"  = 1.13
    "U5"  = 1.67
    "V5"  = 2.1
    "W5"  = 2.38
    "X5"  = 1.53
    "Z5"  = 9.5
    "AA5" = 11
    "AB5" = 21
    "AC5" = 23
    "AE5" = 5.5
    "AJ5" = 7
    "AK5" = 15
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
